# Tokenizer coverage refactor:
#  1. Insert two new paragraphs ("[   ]" and a smart-quote pair "“  “")
#     plus one new empty paragraph, right after the empty paragraph that
#     follows "[Bracket Problem]" and before the "UPPER PROBLEM" paragraph.
#  2. Collapse the "(left punctuation and right punctuation)" paragraph
#     (currently split across three runs with proofErr gramStart/gramEnd
#     markers around the word "left") into a single run with the same text.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: locate the empty paragraph following "[Bracket Problem]".
# ---------------------------------------------------------------------
$bracketIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text.TrimEnd([char]13)
    if ($t -eq "[Bracket Problem]") {
        $bracketIndex = $i
        break
    }
}

$afterBracketEmpty = $d.Paragraphs($bracketIndex + 1)

$insertionRange = $afterBracketEmpty.Range
$insertionRange.Collapse(0)   # wdCollapseEnd

# Create three fresh paragraph marks after the empty paragraph.
$insertionRange.InsertParagraphAfter()
$insertionRange.InsertParagraphAfter()
$insertionRange.InsertParagraphAfter()

# Fill in the text of the first two new paragraphs; the third stays empty.
$d.Paragraphs($bracketIndex + 2).Range.Text = "[   ]"
$d.Paragraphs($bracketIndex + 3).Range.Text = [char]0x201C + "  " + [char]0x201C

# ---------------------------------------------------------------------
# Step 2: find the "(left punctuation and right punctuation)" paragraph
# and collapse its three runs + proofErr markers into a single run.
# ---------------------------------------------------------------------
$targetText = "(left punctuation and right punctuation)"
$punctIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text.TrimEnd([char]13)
    if ($t -eq $targetText) {
        $punctIndex = $i
        break
    }
}

$pr = $d.Paragraphs($punctIndex).Range
[void]$pr.MoveEnd(1, -1)   # wdCharacter — exclude the paragraph mark itself
$pr.Delete()

$pr2 = $d.Paragraphs($punctIndex).Range
[void]$pr2.MoveEnd(1, -1)
$pr2.Text = $targetText
